# Generate Report for Handback
# Updates the localization-status workbook to reflect that the zh-cn and
# de-de handback packages have been generated: status flips from
# "Ready for handoff" to "Handed back: in sync with en-US", and the
# "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns get populated for each language sheet.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: flip the per-language status text -------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

# --- zh-cn sheet -----------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"

$wsZh.Range("I2").Value = "171df3e6-c226-4652-a681-bb125c39f058.md"
$wsZh.Range("J2").Value = "171df3e6-c226-4652-a681-bb125c39f058.0e9125d146a15a4f92304318562cec2d89148c6a.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-21 01:06:10"

$wsZh.Range("I3").Value = "d9077830-64b5-469f-b80c-d17bb6746bb1.md"
$wsZh.Range("J3").Value = "d9077830-64b5-469f-b80c-d17bb6746bb1.43cee252e034940ac71ddd8b71e1a91d37b0fb93.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-21 01:06:10"

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9b47edfe91dd2779a88e5ce69427a492e9740e01/e2e/171df3e6-c226-4652-a681-bb125c39f058.md", "", "", "171df3e6-c226-4652-a681-bb125c39f058.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9b47edfe91dd2779a88e5ce69427a492e9740e01/e2e/d9077830-64b5-469f-b80c-d17bb6746bb1.md", "", "", "d9077830-64b5-469f-b80c-d17bb6746bb1.md")

$wsZh.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsZh.Columns.Item(9).ColumnWidth = 40
$wsZh.Columns.Item(10).ColumnWidth = 40

# --- de-de sheet -------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

$wsDe.Range("I2").Value = "171df3e6-c226-4652-a681-bb125c39f058.md"
$wsDe.Range("J2").Value = "171df3e6-c226-4652-a681-bb125c39f058.0e9125d146a15a4f92304318562cec2d89148c6a.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-21 01:06:17"

$wsDe.Range("I3").Value = "d9077830-64b5-469f-b80c-d17bb6746bb1.md"
$wsDe.Range("J3").Value = "d9077830-64b5-469f-b80c-d17bb6746bb1.43cee252e034940ac71ddd8b71e1a91d37b0fb93.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-21 01:06:17"

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9b47edfe91dd2779a88e5ce69427a492e9740e01/e2e/171df3e6-c226-4652-a681-bb125c39f058.md", "", "", "171df3e6-c226-4652-a681-bb125c39f058.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9b47edfe91dd2779a88e5ce69427a492e9740e01/e2e/d9077830-64b5-469f-b80c-d17bb6746bb1.md", "", "", "d9077830-64b5-469f-b80c-d17bb6746bb1.md")

$wsDe.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsDe.Columns.Item(9).ColumnWidth = 40
$wsDe.Columns.Item(10).ColumnWidth = 40
